$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 (no longer needed) so the used range shrinks to A1:A2
$ws.Rows("3:4").Delete()

# Update the value of A2 for dynamic word selection
$ws.Range("A2").Value = "（我这里用的是【3.11.174】版本"
